# Auto-generated edit script: updates the "想去人数" (F column) values
# on sheets "展览" (sheet 1) and "全部类型" (sheet 4) to match the
# refreshed scrape output (gh-pages commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# --- Sheet 1: 展览 ---
$ws1.Cells.Item(2, 6).Value = 0
$ws1.Cells.Item(4, 6).Value = 19553
$ws1.Cells.Item(5, 6).Value = 783
$ws1.Cells.Item(6, 6).Value = 0
$ws1.Cells.Item(7, 6).Value = 1089
$ws1.Cells.Item(8, 6).Value = 9
$ws1.Cells.Item(9, 6).Value = 7402
$ws1.Cells.Item(10, 6).Value = 484
$ws1.Cells.Item(15, 6).Value = 0
$ws1.Cells.Item(17, 6).Value = 0
$ws1.Cells.Item(18, 6).Value = 0
$ws1.Cells.Item(19, 6).Value = 0
$ws1.Cells.Item(20, 6).Value = 0
$ws1.Cells.Item(22, 6).Value = 0
$ws1.Cells.Item(24, 6).Value = 58
$ws1.Cells.Item(25, 6).Value = 0
$ws1.Cells.Item(26, 6).Value = 1070
$ws1.Cells.Item(28, 6).Value = 8
$ws1.Cells.Item(29, 6).Value = 165
$ws1.Cells.Item(30, 6).Value = 0
$ws1.Cells.Item(31, 6).Value = 553
$ws1.Cells.Item(32, 6).Value = 48
$ws1.Cells.Item(33, 6).Value = 2758
$ws1.Cells.Item(34, 6).Value = 24
$ws1.Cells.Item(35, 6).Value = 0
$ws1.Cells.Item(36, 6).Value = 0
$ws1.Cells.Item(37, 6).Value = 12496
$ws1.Cells.Item(38, 6).Value = 0
$ws1.Cells.Item(39, 6).Value = 58
$ws1.Cells.Item(41, 6).Value = 52
$ws1.Cells.Item(42, 6).Value = 249
$ws1.Cells.Item(44, 6).Value = 0
$ws1.Cells.Item(45, 6).Value = 318

# --- Sheet 4: 全部类型 ---
$ws4.Cells.Item(2, 6).Value = 218
$ws4.Cells.Item(3, 6).Value = 0
$ws4.Cells.Item(4, 6).Value = 19553
$ws4.Cells.Item(5, 6).Value = 0
$ws4.Cells.Item(6, 6).Value = 0
$ws4.Cells.Item(8, 6).Value = 0
$ws4.Cells.Item(9, 6).Value = 7402
$ws4.Cells.Item(10, 6).Value = 484
$ws4.Cells.Item(12, 6).Value = 0
$ws4.Cells.Item(14, 6).Value = 0
$ws4.Cells.Item(15, 6).Value = 102
$ws4.Cells.Item(20, 6).Value = 0
$ws4.Cells.Item(22, 6).Value = 0
$ws4.Cells.Item(23, 6).Value = 48
$ws4.Cells.Item(26, 6).Value = 0
$ws4.Cells.Item(27, 6).Value = 0
$ws4.Cells.Item(28, 6).Value = 0
$ws4.Cells.Item(30, 6).Value = 0
$ws4.Cells.Item(31, 6).Value = 553
$ws4.Cells.Item(32, 6).Value = 0
$ws4.Cells.Item(34, 6).Value = 0
$ws4.Cells.Item(35, 6).Value = 2758
$ws4.Cells.Item(37, 6).Value = 0
$ws4.Cells.Item(39, 6).Value = 12496
$ws4.Cells.Item(40, 6).Value = 1319
$ws4.Cells.Item(41, 6).Value = 58
$ws4.Cells.Item(42, 6).Value = 13
$ws4.Cells.Item(43, 6).Value = 0
$ws4.Cells.Item(44, 6).Value = 0
$ws4.Cells.Item(47, 6).Value = 0

